$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '60.672.53'
Set-TextValue $ws.Range("E2") '  -1.59%  '
Set-TextValue $ws.Range("D3") '2.908.44'
Set-TextValue $ws.Range("E3") '  -1.96%  '
Set-TextValue $ws.Range("E4") '  -0.01%  '
Set-TextValue $ws.Range("D5") '529.75'
Set-TextValue $ws.Range("E5") '  -2.19%  '
Set-TextValue $ws.Range("D6") '144.28'
Set-TextValue $ws.Range("E6") '  -5.37%  '
Set-TextValue $ws.Range("E7") '  -0.01%  '
Set-TextValue $ws.Range("D8") '0.558'
Set-TextValue $ws.Range("E8") '  -0.82%  '
Set-TextValue $ws.Range("D9") '2.916.09'
Set-TextValue $ws.Range("E9") '  -1.95%  '
Set-TextValue $ws.Range("E10") '  -3.27%  '
Set-TextValue $ws.Range("D11") '6.05'
Set-TextValue $ws.Range("E11") '  -0.74%  '
Set-TextValue $ws.Range("E12") '  -1.13%  '
Set-TextValue $ws.Range("D13") '3.418.23'
Set-TextValue $ws.Range("E13") '  -1.91%  '
Set-TextValue $ws.Range("E14") '  +2.29%  '
Set-TextValue $ws.Range("D15") '60.649.22'
Set-TextValue $ws.Range("E15") '  -1.71%  '
Set-TextValue $ws.Range("D16") '22.85'
Set-TextValue $ws.Range("E16") '  -3.45%  '
Set-TextValue $ws.Range("D17") '2.910.45'
Set-TextValue $ws.Range("E17") '  -2.17%  '
Set-TextValue $ws.Range("D18") '0.0000142'
Set-TextValue $ws.Range("E18") '  -3.31%  '
Set-TextValue $ws.Range("E19") '  -1.81%  '
Set-TextValue $ws.Range("D20") '11.74'
Set-TextValue $ws.Range("E20") '  -1.99%  '
Set-TextValue $ws.Range("D21") '363.60'
Set-TextValue $ws.Range("E21") '  -4.45%  '
Set-TextValue $ws.Range("D22") '6.65'
Set-TextValue $ws.Range("E22") '  -0.05%  '
Set-TextValue $ws.Range("E23") '  +0.00%  '
Set-TextValue $ws.Range("E24") '  -0.01%  '
Set-TextValue $ws.Range("D25") '64.94'
Set-TextValue $ws.Range("E25") '  -0.28%  '
Set-TextValue $ws.Range("E26") '  -3.07%  '
Set-TextValue $ws.Range("D27") '0.182'
Set-TextValue $ws.Range("E27") '  -2.48%  '
Set-TextValue $ws.Range("D28") '0.998'
Set-TextValue $ws.Range("E28") '  +0.06%  '
Set-TextValue $ws.Range("D29") '7.91'
Set-TextValue $ws.Range("E29") '  -4.68%  '
Set-TextValue $ws.Range("D30") '0.0₃0865'
Set-TextValue $ws.Range("E30") '  -6.81%  '
Set-TextValue $ws.Range("E31") '  +0.05%  '
Set-TextValue $ws.Range("D32") '1.69'
Set-TextValue $ws.Range("E32") '  -1.95%  '
Set-TextValue $ws.Range("E33") '  -3.07%  '
Set-TextValue $ws.Range("D34") '152.65'
Set-TextValue $ws.Range("E34") '  -3.93%  '
Set-TextValue $ws.Range("D35") '4.40'
Set-TextValue $ws.Range("E35") '  -5.14%  '
Set-TextValue $ws.Range("D36") '5.60'
Set-TextValue $ws.Range("E36") '  -5.60%  '
Set-TextValue $ws.Range("E37") '  -4.83%  '
Set-TextValue $ws.Range("E38") '  -5.03%  '
Set-TextValue $ws.Range("D39") '37.72'
Set-TextValue $ws.Range("E39") '  +1.74%  '
Set-TextValue $ws.Range("E40") '  -3.43%  '
Set-TextValue $ws.Range("D41") '3.74'
Set-TextValue $ws.Range("E41") '  -4.52%  '
Set-TextValue $ws.Range("D42") '2.305.34'
Set-TextValue $ws.Range("E42") '  -4.56%  '
Set-TextValue $ws.Range("D43") '0.649'
Set-TextValue $ws.Range("E43") '  -2.15%  '
Set-TextValue $ws.Range("E44") '  -1.56%  '
Set-TextValue $ws.Range("D45") '20.59'
Set-TextValue $ws.Range("E45") '  -7.37%  '
Set-TextValue $ws.Range("E46") '  -0.03%  '
Set-TextValue $ws.Range("D47") '5.04'
Set-TextValue $ws.Range("E47") '  +1.83%  '
Set-TextValue $ws.Range("E48") '  -2.39%  '
Set-TextValue $ws.Range("B49") 'Stellar'
Set-TextValue $ws.Range("C49") 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D49") '0.0927'
Set-TextValue $ws.Range("E49") '  -2.47%  '
Set-TextValue $ws.Range("B50") 'WhiteBITCoin'
Set-TextValue $ws.Range("C50") 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue $ws.Range("D50") '10.32'
Set-TextValue $ws.Range("E50") '  -1.42%  '
Set-TextValue $ws.Range("D51") '252.56'
Set-TextValue $ws.Range("E51") '  -5.70%  '
